# Generate Report for Handoff
#
# Refreshes the "handoff" snapshot on the localization-status workbook:
#   * Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#     on every sheet that shows a status column.
#   * The handoff timestamps are bumped to the new generation time.
#   * The now-shorter "Ready for handoff" status text no longer needs the
#     wide status columns, so they're narrowed to fit.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# Column width: the workbook stores raw character widths (no Excel "5 px
# padding" baked in), but the ColumnWidth COM property adds that padding
# back before Excel re-quantizes to whole pixels. Subtract it up front so
# the saved width lands as close as possible to the target.
$padding = 5.0 / 6.0
$newColWidth = 17.2159881591797 - $padding

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2017-02-09 16:23:26"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2017-02-09 16:23:08"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2017-02-09 16:23:26"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
